$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.188.90"
$ws.Cells.Item(2, 5).Value = "  -1.32%  "
$ws.Cells.Item(3, 4).Value = "1.658.89"
$ws.Cells.Item(3, 5).Value = "  -0.90%  "
$ws.Cells.Item(4, 5).Value = "  +0.29%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "217.06"
$ws.Cells.Item(5, 5).Value = "  -1.60%  "
$ws.Cells.Item(6, 5).Value = "  -2.13%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2639"
$ws.Cells.Item(8, 5).Value = "  -1.48%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06273"
$ws.Cells.Item(9, 5).Value = "  -1.81%  "
$ws.Cells.Item(10, 5).Value = "  -4.84%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07767"
$ws.Cells.Item(11, 5).Value = "  -0.44%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "4.479"
$ws.Cells.Item(12, 5).Value = "  -0.31%  "
$ws.Cells.Item(13, 4).Value = "1.649.76"
$ws.Cells.Item(13, 5).Value = "  -1.48%  "
$ws.Cells.Item(14, 4).Value = "1.885.70"
$ws.Cells.Item(14, 5).Value = "  -0.91%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.5467"
$ws.Cells.Item(15, 5).Value = "  -1.81%  "
$ws.Cells.Item(16, 4).Value = "0.0₅8130"
$ws.Cells.Item(16, 5).Value = "  -2.52%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "64.96"
$ws.Cells.Item(17, 5).Value = "  -1.06%  "
$ws.Cells.Item(18, 4).Value = "26.200.83"
$ws.Cells.Item(18, 5).Value = "  -1.25%  "
$ws.Cells.Item(19, 5).Value = "  +0.35%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.613"
$ws.Cells.Item(20, 5).Value = "  -3.24%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "192.15"
$ws.Cells.Item(21, 5).Value = "  -0.71%  "
$ws.Cells.Item(22, 5).Value = "  -2.61%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "6.011"
$ws.Cells.Item(23, 5).Value = "  -4.80%  "
$ws.Cells.Item(24, 5).Value = "  +0.35%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "139.41"
$ws.Cells.Item(25, 5).Value = "  -0.10%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.1222"
$ws.Cells.Item(26, 5).Value = "  -4.11%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.281"
$ws.Cells.Item(27, 5).Value = "  -1.76%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "16.15"
$ws.Cells.Item(29, 5).Value = "  +0.75%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.05939"
$ws.Cells.Item(30, 5).Value = "  -4.59%  "
$ws.Cells.Item(31, 5).Value = "  -1.46%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "3.549"
$ws.Cells.Item(32, 5).Value = "  -1.98%  "
$ws.Cells.Item(33, 5).Value = "  -4.28%  "
$ws.Cells.Item(34, 5).Value = "  -6.17%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.9602"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.419"
$ws.Cells.Item(36, 5).Value = "  +0.27%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.767"
$ws.Cells.Item(37, 5).Value = "  -0.50%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.5685"
$ws.Cells.Item(38, 5).Value = "  -6.97%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.042"
$ws.Cells.Item(39, 5).Value = "  -0.31%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.01590"
$ws.Cells.Item(40, 5).Value = "  -1.74%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.8533"
$ws.Cells.Item(41, 5).Value = "  -0.88%  "
$ws.Cells.Item(42, 5).Value = "  +0.32%  "
$ws.Cells.Item(43, 4).Value = "1.012.65"
$ws.Cells.Item(43, 5).Value = "  -7.15%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "101.20"
$ws.Cells.Item(44, 5).Value = "  +0.47%  "
$ws.Cells.Item(45, 4).Value = "1.800.52"
$ws.Cells.Item(45, 5).Value = "  -0.98%  "
$ws.Cells.Item(46, 4).Value = "0.0₈108"
$ws.Cells.Item(46, 5).Value = "  -3.47%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "56.50"
$ws.Cells.Item(47, 5).Value = "  -3.35%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "8.063"
$ws.Cells.Item(49, 5).Value = "  -0.76%  "
$ws.Cells.Item(50, 5).Value = "  -0.56%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.4228"
$ws.Cells.Item(51, 5).Value = "  +0.05%  "
